$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.268.83"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").Value = "3.414.16"
$ws.Range("E3").Value = "  +2.19%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'255.85"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("D6").Value = "'668.67"
$ws.Range("E6").Value = "  +1.44%  "

$ws.Range("D7").Value = "'1.46"
$ws.Range("E7").Value = "  -5.28%  "

$ws.Range("D8").Value = "'0.436"
$ws.Range("E8").Value = "  -4.68%  "

$ws.Range("E9").Value = "  -2.46%  "

$ws.Range("D11").Value = "3.411.59"
$ws.Range("E11").Value = "  +2.22%  "

$ws.Range("E12").Value = "  +3.15%  "

$ws.Range("D13").Value = "'42.05"
$ws.Range("E13").Value = "  -2.23%  "

$ws.Range("D14").Value = "'6.42"
$ws.Range("E14").Value = "  +14.53%  "

$ws.Range("D15").Value = "98.079.96"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").Value = "'0.0000267"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").Value = "4.048.14"
$ws.Range("E17").Value = "  +1.92%  "

$ws.Range("D18").Value = "'8.99"
$ws.Range("E18").Value = "  +19.59%  "

$ws.Range("B19").Value = "Stellar"
$ws.Range("C19").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D19").Value = "'0.589"
$ws.Range("E19").Value = "  +34.93%  "

$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.406.14"
$ws.Range("E20").Value = "  +1.88%  "

$ws.Range("D21").Value = "'17.62"
$ws.Range("E21").Value = "  +4.58%  "

$ws.Range("D22").Value = "'11.04"
$ws.Range("E22").Value = "  +5.08%  "

$ws.Range("D23").Value = "'3.45"
$ws.Range("E23").Value = "  -4.70%  "

$ws.Range("D24").Value = "'511.71"
$ws.Range("E24").Value = "  -3.45%  "

$ws.Range("E25").Value = "  -2.11%  "

$ws.Range("E26").Value = "  +5.95%  "

$ws.Range("D27").Value = "'101.51"
$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("D28").Value = "'12.84"
$ws.Range("E28").Value = "  +1.98%  "

$ws.Range("D29").Value = "3.602.52"
$ws.Range("E29").Value = "  +2.24%  "

$ws.Range("D30").Value = "'0.152"
$ws.Range("E30").Value = "  +1.49%  "

$ws.Range("D31").Value = "'11.56"
$ws.Range("E31").Value = "  +4.84%  "

$ws.Range("E32").Value = "  +3.18%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").Value = "'2.51"
$ws.Range("E34").Value = "  +19.42%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "'0.576"
$ws.Range("E36").Value = "  +7.41%  "

$ws.Range("D37").Value = "'30.05"
$ws.Range("E37").Value = "  +2.52%  "

$ws.Range("D38").Value = "'1.51"
$ws.Range("E38").Value = "  +14.05%  "

$ws.Range("D39").Value = "'7.92"
$ws.Range("E39").Value = "  +1.16%  "

$ws.Range("D40").Value = "'536.25"
$ws.Range("E40").Value = "  +1.41%  "

$ws.Range("E41").Value = "  -3.30%  "

$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").Value = "'0.876"
$ws.Range("E43").Value = "  +6.63%  "

$ws.Range("D44").Value = "'24.71"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "'3.79"
$ws.Range("E45").Value = "  +0.55%  "

$ws.Range("D46").Value = "'8.97"
$ws.Range("E46").Value = "  +13.04%  "

$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'5.83"
$ws.Range("E47").Value = "  +13.61%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0432"
$ws.Range("E48").Value = "  -3.12%  "

$ws.Range("D49").Value = "'1.73"
$ws.Range("E49").Value = "  +15.55%  "

$ws.Range("D50").Value = "'3.27"
$ws.Range("E50").Value = "  -2.68%  "

$ws.Range("D51").Value = "'54.10"
$ws.Range("E51").Value = "  +9.78%  "
